$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.65739733333333
$ws.Range("H2").Value = 37.972192
$ws.Range("I2").Value = 0.2956675086955257
$ws.Range("J2").Value = 0.2956675086955256
$ws.Range("M2").Value = 1.699817666666666
$ws.Range("N2").Value = 5.099453
$ws.Range("O2").Value = 0.748079722752454
$ws.Range("P2").Value = 0.748079722752454
$ws.Range("Q2").Value = 21.51526760121955
$ws.Range("R2").Value = 193.637408410976
$ws.Range("S2").Value = 0.2211828679318576
$ws.Range("T2").Value = 0.2211828679318576
$ws.Range("G3").Value = 12.65739733333333
$ws.Range("H3").Value = 37.972192
$ws.Range("I3").Value = 0.2956675086955257
$ws.Range("J3").Value = 0.2956675086955256
$ws.Range("M3").Value = 0.5724236666666667
$ws.Range("N3").Value = 1.717271
$ws.Range("O3").Value = 0.2519202772475459
$ws.Range("P3").Value = 0.2519202772475459
$ws.Range("Q3").Value = 7.245393792003556
$ws.Range("R3").Value = 65.208544128032
$ws.Range("S3").Value = 0.07448464076366801
$ws.Range("T3").Value = 0.074484640763668
$ws.Range("I4").Value = 0.2489429581834726
$ws.Range("J4").Value = 0.2489429581834726
$ws.Range("M4").Value = 1.699817666666666
$ws.Range("N4").Value = 5.099453
$ws.Range("O4").Value = 0.748079722752454
$ws.Range("P4").Value = 0.748079722752454
$ws.Range("Q4").Value = 18.11519428153411
$ws.Range("R4").Value = 163.036748533807
$ws.Range("S4").Value = 0.1862291791390679
$ws.Range("T4").Value = 0.1862291791390679
$ws.Range("I5").Value = 0.2489429581834726
$ws.Range("J5").Value = 0.2489429581834726
$ws.Range("M5").Value = 0.5724236666666667
$ws.Range("N5").Value = 1.717271
$ws.Range("O5").Value = 0.2519202772475459
$ws.Range("P5").Value = 0.2519202772475459
$ws.Range("Q5").Value = 6.10039896417211
$ws.Range("R5").Value = 54.903590677549
$ws.Range("S5").Value = 0.06271377904440464
$ws.Range("T5").Value = 0.06271377904440464
$ws.Range("G6").Value = 8.42886
$ws.Range("H6").Value = 25.28658
$ws.Range("I6").Value = 0.1968919811642716
$ws.Range("J6").Value = 0.1968919811642716
$ws.Range("M6").Value = 1.699817666666666
$ws.Range("N6").Value = 5.099453
$ws.Range("O6").Value = 0.748079722752454
$ws.Range("P6").Value = 0.748079722752454
$ws.Range("Q6").Value = 14.32752513786
$ws.Range("R6").Value = 128.94772624074
$ws.Range("S6").Value = 0.1472908986815497
$ws.Range("T6").Value = 0.1472908986815497
$ws.Range("G7").Value = 8.42886
$ws.Range("H7").Value = 25.28658
$ws.Range("I7").Value = 0.1968919811642716
$ws.Range("J7").Value = 0.1968919811642716
$ws.Range("M7").Value = 0.5724236666666667
$ws.Range("N7").Value = 1.717271
$ws.Range("O7").Value = 0.2519202772475459
$ws.Range("P7").Value = 0.2519202772475459
$ws.Range("Q7").Value = 4.82487894702
$ws.Range("R7").Value = 43.42391052318
$ws.Range("S7").Value = 0.04960108248272188
$ws.Range("T7").Value = 0.04960108248272188
$ws.Range("G8").Value = 4.628994666666666
$ws.Range("H8").Value = 13.886984
$ws.Range("I8").Value = 0.1081299168237279
$ws.Range("J8").Value = 0.1081299168237279
$ws.Range("M8").Value = 1.699817666666666
$ws.Range("N8").Value = 5.099453
$ws.Range("O8").Value = 0.748079722752454
$ws.Range("P8").Value = 0.748079722752454
$ws.Range("Q8").Value = 7.868446913305776
$ws.Range("R8").Value = 70.81602221975199
$ws.Range("S8").Value = 0.08088979819874026
$ws.Range("T8").Value = 0.08088979819874025
$ws.Range("G9").Value = 4.628994666666666
$ws.Range("H9").Value = 13.886984
$ws.Range("I9").Value = 0.1081299168237279
$ws.Range("J9").Value = 0.1081299168237279
$ws.Range("M9").Value = 0.5724236666666667
$ws.Range("N9").Value = 1.717271
$ws.Range("O9").Value = 0.2519202772475459
$ws.Range("P9").Value = 0.2519202772475459
$ws.Range("Q9").Value = 2.649746100073778
$ws.Range("R9").Value = 23.847714900664
$ws.Range("S9").Value = 0.0272401186249876
$ws.Range("T9").Value = 0.0272401186249876
$ws.Range("G10").Value = 6.437173
$ws.Range("H10").Value = 19.311519
$ws.Range("I10").Value = 0.1503676351330023
$ws.Range("J10").Value = 0.1503676351330023
$ws.Range("M10").Value = 1.699817666666666
$ws.Range("N10").Value = 5.099453
$ws.Range("O10").Value = 0.748079722752454
$ws.Range("P10").Value = 0.748079722752454
$ws.Range("Q10").Value = 10.94202038878967
$ws.Range("R10").Value = 98.47818349910699
$ws.Range("S10").Value = 0.1124869788012385
$ws.Range("T10").Value = 0.1124869788012385
$ws.Range("G11").Value = 6.437173
$ws.Range("H11").Value = 19.311519
$ws.Range("I11").Value = 0.1503676351330023
$ws.Range("J11").Value = 0.1503676351330023
$ws.Range("M11").Value = 0.5724236666666667
$ws.Range("N11").Value = 1.717271
$ws.Range("O11").Value = 0.2519202772475459
$ws.Range("P11").Value = 0.2519202772475459
$ws.Range("Q11").Value = 3.684790171627667
$ws.Range("R11").Value = 33.163111544649
$ws.Range("S11").Value = 0.03788065633176376
$ws.Range("T11").Value = 0.03788065633176375
